$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values for rows 36-43: "Succinate 117/73" -> "Succinate 117/99"
foreach ($r in 36..43) {
    $ws.Cells.Item($r, 3).Value = "Succinate 117/99"
}

# Update the sheet's visible selection (also clears the scrolled topLeftCell
# since C67 becomes the new top-left after selecting/scrolling to it).
$ws.Range("C67").Select()
